# Applies the 2026-01-18 16:37:37 scrape refresh to the 'horarios-141' workbook.
# Updates header metadata (last-updated stamp + row counts) and the arrival-time
# sorted schedule rows on all three sheets (LP1912, LP1912-215, 6203-6173),
# including newly scraped rows appended at the end of each table.

$wb = $excel.ActiveWorkbook

# ===== Sheet: LP1912 =====
$ws = $wb.Worksheets.Item('LP1912')
$ws.Cells.Item(2,1).Value = 'Última actualización: 16:37:37'
$ws.Cells.Item(3,1).Value = 'Total filas: 259'
$ws.Cells.Item(20,1).Value = '06:16:41'
$ws.Cells.Item(20,3).Value = '23_HERNANDEZ'
$ws.Cells.Item(20,4).Value = 65
$ws.Cells.Item(21,1).Value = '05:57:04'
$ws.Cells.Item(21,3).Value = '16_SANTA ANA'
$ws.Cells.Item(21,4).Value = 84
$ws.Cells.Item(86,1).Value = '10:04:30'
$ws.Cells.Item(86,3).Value = '14_ABASTO'
$ws.Cells.Item(86,4).Value = 25
$ws.Cells.Item(87,1).Value = '08:38:24'
$ws.Cells.Item(87,3).Value = '15_ABASTO'
$ws.Cells.Item(87,4).Value = 111
$ws.Cells.Item(102,1).Value = '09:22:34'
$ws.Cells.Item(102,3).Value = '15X38_ABASTO'
$ws.Cells.Item(102,4).Value = 113
$ws.Cells.Item(103,1).Value = '10:56:15'
$ws.Cells.Item(103,3).Value = '14_ABASTO'
$ws.Cells.Item(103,4).Value = 19
$ws.Cells.Item(118,1).Value = '11:33:52'
$ws.Cells.Item(118,3).Value = '23_HERNANDEZ'
$ws.Cells.Item(118,4).Value = 19
$ws.Cells.Item(119,1).Value = '10:36:50'
$ws.Cells.Item(119,3).Value = '225_GOMEZ'
$ws.Cells.Item(119,4).Value = 76
$ws.Cells.Item(133,1).Value = '11:33:52'
$ws.Cells.Item(133,3).Value = '23_HERNANDEZ'
$ws.Cells.Item(133,4).Value = 59
$ws.Cells.Item(134,1).Value = '10:36:50'
$ws.Cells.Item(134,3).Value = '14_ABASTO'
$ws.Cells.Item(134,4).Value = 116
$ws.Cells.Item(135,3).Value = '27_EL RETIRO'
$ws.Cells.Item(136,3).Value = '14_ABASTO'
$ws.Cells.Item(184,1).Value = '14:32:44'
$ws.Cells.Item(184,3).Value = '14X44_ABASTO'
$ws.Cells.Item(184,4).Value = 1
$ws.Cells.Item(185,1).Value = '13:55:43'
$ws.Cells.Item(185,3).Value = '215C_EL PATO'
$ws.Cells.Item(185,4).Value = 38
$ws.Cells.Item(253,1).Value = '16:37:37'
$ws.Cells.Item(253,2).Value = '17:53'
$ws.Cells.Item(253,3).Value = '23_HERNANDEZ'
$ws.Cells.Item(253,4).Value = 76
$ws.Cells.Item(254,1).Value = '16:12:06'
$ws.Cells.Item(254,2).Value = '17:58'
$ws.Cells.Item(254,3).Value = '17_ROMERO'
$ws.Cells.Item(254,4).Value = 106
$ws.Cells.Item(255,1).Value = '16:28:21'
$ws.Cells.Item(255,2).Value = '18:05'
$ws.Cells.Item(255,4).Value = 97
$ws.Cells.Item(256,2).Value = '18:06'
$ws.Cells.Item(256,3).Value = '11_ETCHEVERRY'
$ws.Cells.Item(256,4).Value = 114
$ws.Cells.Item(257,3).Value = '16_P MOR-SANTA ANA'
$ws.Cells.Item(258,1).Value = '16:12:06'
$ws.Cells.Item(258,2).Value = '18:10'
$ws.Cells.Item(258,3).Value = '15_ABASTO'
$ws.Cells.Item(258,4).Value = 118
$ws.Cells.Item(259,2).Value = '18:17'
$ws.Cells.Item(259,3).Value = '10_OLMOS'
$ws.Cells.Item(259,4).Value = 109
$ws.Cells.Item(260,1).Value = '16:37:37'
$ws.Cells.Item(260,2).Value = '18:21'
$ws.Cells.Item(260,3).Value = '215C_EL PATO'
$ws.Cells.Item(260,4).Value = 104
$ws.Cells.Item(261,1).Value = '16:28:21'
$ws.Cells.Item(261,2).Value = '18:22'
$ws.Cells.Item(261,3).Value = '215C_EL PATO'
$ws.Cells.Item(261,4).Value = 114
$ws.Cells.Item(261,5).Value = 'LP1912'
$ws.Cells.Item(262,1).Value = '16:28:21'
$ws.Cells.Item(262,2).Value = '18:25'
$ws.Cells.Item(262,3).Value = '16_SANTA ANA'
$ws.Cells.Item(262,4).Value = 117
$ws.Cells.Item(262,5).Value = 'LP1912'
$ws.Cells.Item(263,1).Value = '16:37:37'
$ws.Cells.Item(263,2).Value = '18:29'
$ws.Cells.Item(263,3).Value = '14_ABASTO'
$ws.Cells.Item(263,4).Value = 112
$ws.Cells.Item(263,5).Value = 'LP1912'
$ws.Cells.Item(264,1).Value = '16:37:37'
$ws.Cells.Item(264,2).Value = '18:36'
$ws.Cells.Item(264,3).Value = '15X38_ABASTO'
$ws.Cells.Item(264,4).Value = 119
$ws.Cells.Item(264,5).Value = 'LP1912'

# ===== Sheet: LP1912-215 =====
$ws = $wb.Worksheets.Item('LP1912-215')
$ws.Cells.Item(2,1).Value = 'Última actualización: 16:37:37'
$ws.Cells.Item(3,1).Value = 'Total filas: 44'
$ws.Cells.Item(48,1).Value = '16:37:37'
$ws.Cells.Item(48,2).Value = '18:21'
$ws.Cells.Item(48,4).Value = 104
$ws.Cells.Item(49,1).Value = '16:28:21'
$ws.Cells.Item(49,2).Value = '18:22'
$ws.Cells.Item(49,3).Value = '215C_EL PATO'
$ws.Cells.Item(49,4).Value = 114
$ws.Cells.Item(49,5).Value = 'LP1912'

# ===== Sheet: 6203-6173 =====
$ws = $wb.Worksheets.Item('6203-6173')
$ws.Cells.Item(2,1).Value = 'Última actualización: 16:37:37'
$ws.Cells.Item(3,1).Value = 'Total filas: 40'
$ws.Cells.Item(45,1).Value = '16:37:37'
$ws.Cells.Item(45,2).Value = '18:35'
$ws.Cells.Item(45,3).Value = '215C_LA PLATA'
$ws.Cells.Item(45,4).Value = 118
$ws.Cells.Item(45,5).Value = 'L6203'

